$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1622.2354
$ws.Range("I15").Value = 1622.2354
$ws.Range("K15").Value = 4866.706200000001
$ws.Range("M15").Value = -4697.706200000001

# Row 17
$ws.Range("H17").Value = 13896.986
$ws.Range("J17").Value = 13896.986
$ws.Range("L17").Value = 41690.958
$ws.Range("N17").Value = -42026.958

# Row 51
$ws.Range("H51").Value = 7833.2856
$ws.Range("J51").Value = 5736.8423
$ws.Range("L51").Value = 5736.8423
$ws.Range("N51").Value = -6704.8423

# Row 113
$ws.Range("H113").Value = 7401.5
$ws.Range("I113").Value = 7800
$ws.Range("J113").Value = 7003
$ws.Range("K113").Value = 7800
$ws.Range("L113").Value = 7003
$ws.Range("N113").Value = -13511
$ws.Range("M113").Value = -4546

# Row 116
$ws.Range("H116").Value = 11987.863
$ws.Range("I116").Value = 7412.4287
$ws.Range("K116").Value = 7412.4287
$ws.Range("M116").Value = -3970.4287

# Row 137
$ws.Range("H137").Value = 1307.5
$ws.Range("I137").Value = 725.2857
$ws.Range("J137").Value = 2666
$ws.Range("K137").Value = 2175.8571
$ws.Range("L137").Value = 7998
$ws.Range("M137").Value = 374.1428999999998
$ws.Range("N137").Value = -13098

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1979.6666
$ws.Range("I45").Value = 1386.1428
$ws.Range("K45").Value = 1386.1428
$ws.Range("M45").Value = -1009.1428

# Row 46
$ws.Range("H46").Value = 21051.273
$ws.Range("I46").Value = 9898
$ws.Range("J46").Value = 22166.6
$ws.Range("K46").Value = 9898
$ws.Range("L46").Value = 22166.6
$ws.Range("M46").Value = -9579
$ws.Range("N46").Value = -22804.6

# Row 61
$ws.Range("H61").Value = 8405537
$ws.Range("I61").Value = 11768132
$ws.Range("J61").Value = 1260024.4
$ws.Range("K61").Value = 11768132
$ws.Range("L61").Value = 1260024.4
$ws.Range("M61").Value = -11767920
$ws.Range("N61").Value = -1260448.4

# Row 74
$ws.Range("H74").Value = 1537.2122
$ws.Range("I74").Value = 1048.7931
$ws.Range("K74").Value = 1048.7931
$ws.Range("M74").Value = -174.7931000000001

# Row 77
$ws.Range("H77").Value = 1537.2122
$ws.Range("I77").Value = 1048.7931
$ws.Range("K77").Value = 5243.9655
$ws.Range("M77").Value = -875.9655000000002

# Row 97
$ws.Range("H97").Value = 846.27905
$ws.Range("I97").Value = 846.27905
$ws.Range("K97").Value = 846.27905
$ws.Range("M97").Value = -350.27905

# Row 110
$ws.Range("H110").Value = 5878.25
$ws.Range("I110").Value = 5465.5386
$ws.Range("K110").Value = 5465.5386
$ws.Range("M110").Value = -3420.5386

# Row 122
$ws.Range("H122").Value = 4352
$ws.Range("J122").Value = 4185.8
$ws.Range("L122").Value = 12557.4
$ws.Range("N122").Value = -17457.4

# Row 132
$ws.Range("H132").Value = 2442746.5
$ws.Range("I132").Value = 3511.8572
$ws.Range("K132").Value = 10535.5716
$ws.Range("M132").Value = -8005.571599999999

# Row 136
$ws.Range("H136").Value = 8405537
$ws.Range("I136").Value = 11768132
$ws.Range("J136").Value = 1260024.4
$ws.Range("K136").Value = 35304396
$ws.Range("L136").Value = 3780073.2
$ws.Range("M136").Value = -35301846
$ws.Range("N136").Value = -3785173.2

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2261.7878
$ws.Range("I94").Value = 1997.7407
$ws.Range("K94").Value = 1997.7407
$ws.Range("M94").Value = -1546.7407

# Row 102
$ws.Range("H102").Value = 24998
$ws.Range("I102").Value = 6247.75
$ws.Range("K102").Value = 6247.75
$ws.Range("M102").Value = -3002.75

# Row 105
$ws.Range("H105").Value = 1528760.5
$ws.Range("I105").Value = 2542967.5
$ws.Range("J105").Value = 7449.8335
$ws.Range("K105").Value = 2542967.5
$ws.Range("L105").Value = 7449.8335
$ws.Range("M105").Value = -2541220.5
$ws.Range("N105").Value = -10943.8335

# Row 107
$ws.Range("H107").Value = 5233
$ws.Range("I107").Value = 5085.8335
$ws.Range("K107").Value = 5085.8335
$ws.Range("M107").Value = -3165.8335

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 11117877
$ws.Range("I16").Value = 20004796
$ws.Range("K16").Value = 20004796
$ws.Range("M16").Value = -20004509

# Row 28
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Row 107
$ws.Range("H107").Value = 4512.5386
$ws.Range("I107").Value = 4100
$ws.Range("J107").Value = 4993.8335
$ws.Range("K107").Value = 4100
$ws.Range("L107").Value = 4993.8335
$ws.Range("M107").Value = -2180
$ws.Range("N107").Value = -8833.833500000001

# Row 113
$ws.Range("H113").Value = 11117877
$ws.Range("I113").Value = 20004796
$ws.Range("K113").Value = 20004796
$ws.Range("M113").Value = -20002626

# Row 122
$ws.Range("H122").Value = 3572.6
$ws.Range("I122").Value = 3310.4
$ws.Range("J122").Value = 4097
$ws.Range("K122").Value = 9931.200000000001
$ws.Range("L122").Value = 12291
$ws.Range("M122").Value = -7481.200000000001
$ws.Range("N122").Value = -17191

# Row 134
$ws.Range("H134").Value = 2326.3684
$ws.Range("I134").Value = 2100.0588
$ws.Range("K134").Value = 6300.176399999999
$ws.Range("M134").Value = -3765.176399999999

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 6521999.5
$ws.Range("I4").Value = 6521999.5
$ws.Range("K4").Value = 19565998.5
$ws.Range("M4").Value = -19565886.5

# Row 5
$ws.Range("H5").Value = 1482.4546
$ws.Range("I5").Value = 700.8889
$ws.Range("K5").Value = 2102.6667
$ws.Range("M5").Value = -1990.6667

# Row 133
$ws.Range("H133").Value = 21664
$ws.Range("I133").Value = 15511.091
$ws.Range("K133").Value = 46533.273
$ws.Range("M133").Value = -41473.273

# Row 134
$ws.Range("H134").Value = 18417.902
$ws.Range("I134").Value = 2509.4666
$ws.Range("K134").Value = 7528.399800000001
$ws.Range("M134").Value = -2458.399800000001

# Row 135
$ws.Range("H135").Value = 1482.4546
$ws.Range("I135").Value = 700.8889
$ws.Range("K135").Value = 6308.0001
$ws.Range("M135").Value = -3773.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 549.62164
$ws.Range("I97").Value = 425.5
$ws.Range("J97").Value = 935.7778
$ws.Range("K97").Value = 425.5
$ws.Range("L97").Value = 935.7778
$ws.Range("M97").Value = 70.5
$ws.Range("N97").Value = -1927.7778

# Row 102
$ws.Range("H102").Value = 2882
$ws.Range("I102").Value = 2880.2222
$ws.Range("K102").Value = 2880.2222
$ws.Range("M102").Value = -1258.2222

# Row 113
$ws.Range("H113").Value = 2648991.5
$ws.Range("I113").Value = 3884.4
$ws.Range("K113").Value = 3884.4
$ws.Range("M113").Value = -1714.4

# Row 122
$ws.Range("H122").Value = 3531.5881
$ws.Range("I122").Value = 3402.4666
$ws.Range("K122").Value = 10207.3998
$ws.Range("M122").Value = -7757.399800000001

# Row 132
$ws.Range("H132").Value = 12502735
$ws.Range("I132").Value = 2720.75
$ws.Range("J132").Value = 25002750
$ws.Range("K132").Value = 8162.25
$ws.Range("L132").Value = 75008250
$ws.Range("M132").Value = -5632.25
$ws.Range("N132").Value = -75013310

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 5950
$ws.Range("I46").Value = 1900
$ws.Range("K46").Value = 1900
$ws.Range("M46").Value = -1712

# Row 93
$ws.Range("H93").Value = 1636405.4
$ws.Range("J93").Value = 11122911
$ws.Range("L93").Value = 11122911
$ws.Range("N93").Value = -11125407

# Row 132
$ws.Range("H132").Value = 2902
$ws.Range("I132").Value = 1882.52
$ws.Range("K132").Value = 5647.559999999999
$ws.Range("M132").Value = -3117.559999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1771
$ws.Range("I81").Value = 1771
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3542
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2481
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 1771
$ws.Range("I84").Value = 1771
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 17710
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -12406
$ws.Range("N84").ClearContents()

# Row 122
$ws.Range("H122").Value = 3233.9048
$ws.Range("I122").Value = 2686.5
$ws.Range("J122").Value = 3731.5454
$ws.Range("K122").Value = 8059.5
$ws.Range("L122").Value = 11194.6362
$ws.Range("M122").Value = -5609.5
$ws.Range("N122").Value = -16094.6362

# Row 132
$ws.Range("H132").Value = 402229.16
$ws.Range("I132").Value = 2265.0557
$ws.Range("K132").Value = 6795.1671
$ws.Range("M132").Value = -4265.1671

# Row 133
$ws.Range("H133").Value = 59465
$ws.Range("J133").Value = 59465
$ws.Range("L133").Value = 59465
$ws.Range("N133").Value = -69585
